$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 647
$ws.Range("F5").Value = 29
$ws.Range("F7").Value = 2160
$ws.Range("G7").Value = 78
$ws.Range("G8").Value = 65
$ws.Range("F9").Value = 853
$ws.Range("G9").Value = 78
$ws.Range("F11").Value = 92
$ws.Range("F13").Value = 324
$ws.Range("F14").Value = 101
$ws.Range("F15").Value = 898
$ws.Range("F18").Value = 1788
$ws.Range("F22").Value = 60
$ws.Range("F24").Value = 1461
$ws.Range("F27").Value = 355
$ws.Range("F28").Value = 625
$ws.Range("F29").Value = 424
$ws.Range("F30").Value = 2473
$ws.Range("F31").Value = 384
$ws.Range("F32").Value = 98
$ws.Range("F35").Value = 480
$ws.Range("F36").Value = 199
$ws.Range("F37").Value = 926
$ws.Range("F38").Value = 708
$ws.Range("F39").Value = 38
$ws.Range("F40").Value = 504
$ws.Range("F41").Value = 486

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 56
$ws.Range("F16").Value = 88
$ws.Range("F22").Value = 122
$ws.Range("F24").Value = 436

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 326

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 647
$ws.Range("F8").Value = 29
$ws.Range("G10").Value = 65
$ws.Range("F11").Value = 853
$ws.Range("G11").Value = 78
$ws.Range("F13").Value = 92
$ws.Range("F14").Value = 324
$ws.Range("F15").Value = 101
$ws.Range("F17").Value = 898
$ws.Range("F21").Value = 326
$ws.Range("F22").Value = 1788
$ws.Range("F26").Value = 56
$ws.Range("F30").Value = 1461
$ws.Range("F34").Value = 355
$ws.Range("F35").Value = 625
$ws.Range("F36").Value = 424
$ws.Range("F37").Value = 98
$ws.Range("F39").Value = 480
$ws.Range("F40").Value = 200
$ws.Range("F41").Value = 926
$ws.Range("F45").Value = 436
$ws.Range("F46").Value = 708
$ws.Range("F47").Value = 38
$ws.Range("F48").Value = 504
$ws.Range("F49").Value = 486
